$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell data for rows 2-7 (columns A:T), replacing the previous rows 2-4.
# Column order: A..T
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..T numeric metrics
$data = @{
    2 = @("ECs", "Sele", "Cd44", "ECs", 3, 1, 7.879565666666667, 23.638697, 0.9977172793687663, 0.9977172793687664, 3, 1, 297.8183156666666, 893.4549469999999, 0.8852156413092672, 0.8852156413092673, 2346.678975031562, 21120.11077528406, 0.8831949413017598, 0.88319494130176)
    3 = @("ECs", "Sele", "Cd44", "FAPs", 3, 1, 7.879565666666667, 23.638697, 0.9977172793687663, 0.9977172793687664, 3, 1, 24.34034433333333, 73.021033, 0.07234764413494278, 0.0723476441349428, 191.7913415237779, 1726.122073714001, 0.07218249467505479, 0.07218249467505482)
    4 = @("ECs", "Sele", "Cd44", "sCs", 3, 1, 7.879565666666667, 23.638697, 0.9977172793687663, 0.9977172793687664, 3, 1, 14.277234, 42.831702, 0.04243671455578994, 0.04243671455578994, 112.498402841366, 1012.485625572294, 0.04233984339195166, 0.04233984339195167)
    5 = @("sCs", "Sele", "Cd44", "ECs", 1, 0.3333333333333333, 0.018028, 0.054084, 0.002282720631233623, 0.002282720631233623, 3, 1, 297.8183156666666, 893.4549469999999, 0.8852156413092672, 0.8852156413092673, 5.369068594838666, 48.321617353548, 0.002020700007507367, 0.002020700007507367)
    6 = @("sCs", "Sele", "Cd44", "FAPs", 1, 0.3333333333333333, 0.018028, 0.054084, 0.002282720631233623, 0.002282720631233623, 3, 1, 24.34034433333333, 73.021033, 0.07234764413494278, 0.0723476441349428, 0.4388077276413333, 3.949269548772, 0.0001651494598879821, 0.0001651494598879822)
    7 = @("sCs", "Sele", "Cd44", "sCs", 1, 0.3333333333333333, 0.018028, 0.054084, 0.002282720631233623, 0.002282720631233623, 3, 1, 14.277234, 42.831702, 0.04243671455578994, 0.04243671455578994, 0.257389974552, 2.316509770968, 0.0000968711638382739, 0.00009687116383827391)
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in 2..7) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowVals[$i]
    }
}
